$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1211596285783969
$ws.Range("D2").Value = 0.17341228526778
$ws.Range("E2").Value = 0.1958613402183147
$ws.Range("F2").Value = 2.973650978365868
$ws.Range("G2").Value = 0.002581764597226413
$ws.Range("I2").Value = 2.431055493695368
$ws.Range("L2").Value = 0.3751323624953784
$ws.Range("M2").Value = 11.40542490509586
$ws.Range("C3").Value = 0.1063735221813147
$ws.Range("D3").Value = 0.1855300351781324
$ws.Range("E3").Value = 0.1724891289907404
$ws.Range("F3").Value = 3.033923446151277
$ws.Range("G3").Value = 0.002596801957842243
$ws.Range("I3").Value = 2.513062227762603
$ws.Range("L3").Value = 0.3303517992685272
$ws.Range("M3").Value = 10.15005367726991
$ws.Range("C4").Value = 0.09728644906013528
$ws.Range("D4").Value = 0.1935263628975861
$ws.Range("E4").Value = 0.1581723004550639
$ws.Range("F4").Value = 3.077850035676036
$ws.Range("G4").Value = 0.002606409777363483
$ws.Range("I4").Value = 2.569372736765075
$ws.Range("L4").Value = 0.3030010355862771
$ws.Range("M4").Value = 9.378925818445964
$ws.Range("C5").Value = 0.09358044754100092
$ws.Range("D5").Value = 0.1969213881353014
$ws.Range("E5").Value = 0.1523449093145288
$ws.Range("F5").Value = 3.097446916794922
$ws.Range("G5").Value = 0.002610420226121106
$ws.Range("I5").Value = 2.593782646694891
$ws.Range("L5").Value = 0.2918873618153555
$ws.Range("M5").Value = 9.064519221373985
$ws.Range("C6").Value = 0.09296486719432551
$ws.Range("D6").Value = 0.1974932753755105
$ws.Range("E6").Value = 0.1513776420958806
$ws.Range("F6").Value = 3.100802180968742
$ws.Range("G6").Value = 0.002611091933801545
$ws.Range("I6").Value = 2.597923239642142
$ws.Range("L6").Value = 0.2900437545594343
$ws.Range("M6").Value = 9.012299908315981
$ws.Range("C7").Value = 0.09723648163401322
$ws.Range("D7").Value = 0.1935716012269744
$ws.Range("E7").Value = 0.1580936846566985
$ws.Range("F7").Value = 3.078107513597033
$ws.Range("G7").Value = 0.00260646347707098
$ws.Range("I7").Value = 2.569696059825063
$ws.Range("L7").Value = 0.3028510287541906
$ws.Range("M7").Value = 9.374686398398126
$ws.Range("C8").Value = 0.1160625799276431
$ws.Range("D8").Value = 0.1774730424506572
$ws.Range("E8").Value = 0.1877945315215754
$ws.Range("F8").Value = 2.992970326919462
$ws.Range("G8").Value = 0.002586872258807113
$ws.Range("I8").Value = 2.45807325760876
$ws.Range("L8").Value = 0.3596592286479279
$ws.Range("M8").Value = 10.97258728525725
$ws.Range("C9").Value = 0.1529569503083792
$ws.Range("D9").Value = 0.1504710758439103
$ws.Range("E9").Value = 0.2463913182318578
$ws.Range("F9").Value = 2.882939245116233
$ws.Range("G9").Value = 0.002551383675795679
$ws.Range("I9").Value = 2.288096441087376
$ws.Range("L9").Value = 0.4724293527338261
$ws.Range("M9").Value = 14.10777249552507
$ws.Range("C10").Value = 0.180114723872606
$ws.Range("D10").Value = 0.1336412050860289
$ws.Range("E10").Value = 0.2897915123188994
$ws.Range("F10").Value = 2.839714507182833
$ws.Range("G10").Value = 0.002527031496223575
$ws.Range("I10").Value = 2.195393085004895
$ws.Range("L10").Value = 0.5564588087297864
$ws.Range("M10").Value = 16.41892212857704
$ws.Range("C11").Value = 0.1924965132400587
$ws.Range("D11").Value = 0.1266901633714355
$ws.Range("E11").Value = 0.3096441272773944
$ws.Range("F11").Value = 2.828891032560819
$ws.Range("G11").Value = 0.00251631189836643
$ws.Range("I11").Value = 2.16075194204484
$ws.Range("L11").Value = 0.595025503728067
$ws.Range("M11").Value = 17.47364262253376
$ws.Range("C12").Value = 0.1971905010301782
$ws.Range("D12").Value = 0.1241639806812387
$ws.Range("E12").Value = 0.317180475682548
$ws.Range("F12").Value = 2.826122660394788
$ws.Range("G12").Value = 0.002512302926399024
$ws.Range("I12").Value = 2.148765026706513
$ws.Range("L12").Value = 0.6096862445621127
$ws.Range("M12").Value = 17.87366625310892
$ws.Range("C13").Value = 0.1961793104578646
$ws.Range("D13").Value = 0.1247032468490374
$ws.Range("E13").Value = 0.3155565136821963
$ws.Range("F13").Value = 2.826658741322802
$ws.Range("G13").Value = 0.00251316411298635
$ws.Range("I13").Value = 2.15129553123603
$ws.Range("L13").Value = 0.6065261620677518
$ws.Range("M13").Value = 17.78748391329754
$ws.Range("C14").Value = 0.1928825769476816
$ws.Range("D14").Value = 0.126480180282158
$ws.Range("E14").Value = 0.3102637575071299
$ws.Range("F14").Value = 2.828636267174602
$ws.Range("G14").Value = 0.002515981076643488
$ws.Range("I14").Value = 2.159742820982601
$ws.Range("L14").Value = 0.5962304780404111
$ws.Range("M14").Value = 17.50653946318283
$ws.Range("C15").Value = 0.1908639595297075
$ws.Range("D15").Value = 0.1275825536595718
$ws.Range("E15").Value = 0.30702430095252
$ws.Range("F15").Value = 2.830022603158085
$ws.Range("G15").Value = 0.002517713062635722
$ws.Range("I15").Value = 2.165065782424676
$ws.Range("L15").Value = 0.589931642062254
$ws.Range("M15").Value = 17.33453846249574
$ws.Range("C16").Value = 0.1793062094352535
$ws.Range("D16").Value = 0.1341100494352787
$ws.Range("E16").Value = 0.2884965399206578
$ws.Range("F16").Value = 2.840604706328946
$ws.Range("G16").Value = 0.00252773915474935
$ws.Range("I16").Value = 2.197812545523618
$ws.Range("L16").Value = 0.5539458620489199
$ws.Range("M16").Value = 16.3500731005625
$ws.Range("C17").Value = 0.1722238693103009
$ws.Range("D17").Value = 0.1382982425919437
$ws.Range("E17").Value = 0.2771603796829396
$ws.Range("F17").Value = 2.849401633014452
$ws.Range("G17").Value = 0.002533980765314414
$ws.Range("I17").Value = 2.219863708073987
$ws.Range("L17").Value = 0.5319621663104783
$ws.Range("M17").Value = 15.74709370284637
$ws.Range("C18").Value = 0.1681528070776039
$ws.Range("D18").Value = 0.1407732403086754
$ws.Range("E18").Value = 0.2706502576705532
$ws.Range("F18").Value = 2.855289709981406
$ws.Range("G18").Value = 0.002537604578664688
$ws.Range("I18").Value = 2.233252695480601
$ws.Range("L18").Value = 0.5193493767601183
$ws.Range("M18").Value = 15.40058378096063
$ws.Range("C19").Value = 0.166774808559893
$ws.Range("D19").Value = 0.1416224402009263
$ws.Range("E19").Value = 0.2684477010445505
$ws.Range("F19").Value = 2.857423824709144
$ws.Range("G19").Value = 0.002538837381787357
$ws.Range("I19").Value = 2.237905763010346
$ws.Range("L19").Value = 0.5150841262566246
$ws.Range("M19").Value = 15.28331028039781
$ws.Range("C20").Value = 0.1729775273033454
$ws.Range("D20").Value = 0.1378455324516636
$ws.Range("E20").Value = 0.2783660643707861
$ws.Range("F20").Value = 2.84837903908479
$ws.Range("G20").Value = 0.002533312845012015
$ws.Range("I20").Value = 2.217442942043618
$ws.Range("L20").Value = 0.5342990415406348
$ws.Range("M20").Value = 15.81124904863577
$ws.Range("C21").Value = 0.1938507535137148
$ws.Range("D21").Value = 0.1259553359408798
$ws.Range("E21").Value = 0.3118178399497111
$ws.Range("F21").Value = 2.828018830911702
$ws.Range("G21").Value = 0.002515152311089992
$ws.Range("I21").Value = 2.15723055782469
$ws.Range("L21").Value = 0.5992529790289893
$ws.Range("M21").Value = 17.58904155658672
$ws.Range("C22").Value = 0.2075240217535281
$ws.Range("D22").Value = 0.1188047586989782
$ws.Range("E22").Value = 0.3337903735673251
$ws.Range("F22").Value = 2.822495648958125
$ws.Range("G22").Value = 0.002503575969158336
$ws.Range("I22").Value = 2.124494435243619
$ws.Range("L22").Value = 0.6420365209879151
$ws.Range("M22").Value = 18.75464124532937
$ws.Range("C23").Value = 0.200223011102338
$ws.Range("D23").Value = 0.1225627862583352
$ws.Range("E23").Value = 0.3220521696799921
$ws.Range("F23").Value = 2.82471093805276
$ws.Range("G23").Value = 0.00250972812538067
$ws.Range("I23").Value = 2.141344366277309
$ws.Range("L23").Value = 0.6191691583765362
$ws.Range("M23").Value = 18.13215016922231
$ws.Range("C24").Value = 0.1726367965369207
$ws.Range("D24").Value = 0.1380499938991875
$ws.Range("E24").Value = 0.2778209521467261
$ws.Range("F24").Value = 2.848838771080636
$ws.Range("G24").Value = 0.002533614701277252
$ws.Range("I24").Value = 2.218535157873376
$ws.Range("L24").Value = 0.5332424600276227
$ws.Range("M24").Value = 15.78224395700471
$ws.Range("C25").Value = 0.1429720389769784
$ws.Range("D25").Value = 0.1572643941632847
$ws.Range("E25").Value = 0.2304883852180666
$ws.Range("F25").Value = 2.906318742587388
$ws.Range("G25").Value = 0.002560677150247992
$ws.Range("I25").Value = 2.328605933953483
$ws.Range("L25").Value = 0.441741799330174
$ws.Range("M25").Value = 13.25878177220386
